# Generate Report for Handoff
#
# The "6e7ac35b-8f2d-4cd6-ace4-af348147644c.md" file has finished its
# handoff-xliff generation step: its status flips from "In Translation"
# to "Ready for handoff" (with a fresh handoff timestamp) on the Overview
# sheet and on each per-language sheet (zh-cn, de-de). The zh-cn/de-de
# rows also pick up the "mt" priority and (for zh-cn) a new "Latest
# Handoff Datetime".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet — row for 6e7ac35b...md (row 3)
# Columns: E = zh-cn status, F = de-de status, G = Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-04 22:16:52"
$wsOverview.Range("E:F").EntireColumn.AutoFit()

# ---------------------------------------------------------------------
# zh-cn sheet — row for 6e7ac35b...md (row 3)
# Columns: C = Status, E = Priority, H = Latest Handoff Datetime
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-09-04 22:16:47"
$wsZhCn.Range("C:C").EntireColumn.AutoFit()

# ---------------------------------------------------------------------
# de-de sheet — row for 6e7ac35b...md (row 3)
# Columns: C = Status, E = Priority, H = Latest Handoff Datetime
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-09-04 22:16:52"
$wsDeDe.Range("C:C").EntireColumn.AutoFit()
